# Generate Report for Handback
#
# The localization-status report is updated to reflect that the
# handback (de-de docs now in sync with en-US, zh-cn's handback
# timestamp refreshed) has completed:
#   - Overview/zh-cn/de-de "Status" columns move from
#     "Ready for handoff" -> "Handed back: in sync with en-US"
#   - zh-cn & de-de sheets get their "Latest Target File" (I) and
#     "Latest Handback File" (J) columns populated with a hyperlink to
#     the source doc / the generated handback xliff file name
#   - "Latest Handback DateTime" (K) is stamped
#   - The two narrow columns get widened so the new values are legible

$wb = $excel.ActiveWorkbook

$srcDocUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/32bbfa4b85de1cf3676f30cc08d7fd5e29ea698c/e2e/291c5604-b16a-4a9d-9f37-928989d8f104.md"
$srcDocName = "291c5604-b16a-4a9d-9f37-928989d8f104.md"

# ---------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------
$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = $statusNew
$ovw.Range("F2").Value = $statusNew
$ovw.Range("E3").Value = $statusNew
$ovw.Range("F3").Value = $statusNew

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusNew
$zhcn.Range("C3").Value = $statusNew

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusNew
$dede.Range("C3").Value = $statusNew

# ---------------------------------------------------------------
# 2. zh-cn: Latest Target File (I), Latest Handback File (J),
#    Latest Handback DateTime (K)
# ---------------------------------------------------------------
$zhcn.Range("J2").Value = "291c5604-b16a-4a9d-9f37-928989d8f104.47401c9d65c9b1103fca9549ff81d630b6a91ca0.zh-cn.xlf"
$zhcn.Range("J3").Value = "291c5604-b16a-4a9d-9f37-928989d8f104.47401c9d65c9b1103fca9549ff81d630b6a91ca0.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-24 15:14:04"
$zhcn.Range("K3").Value = "2016-08-24 15:14:04"

# ---------------------------------------------------------------
# 3. de-de: Latest Target File (I), Latest Handback File (J),
#    Latest Handback DateTime (K)
# ---------------------------------------------------------------
$dede.Range("J2").Value = "291c5604-b16a-4a9d-9f37-928989d8f104.47401c9d65c9b1103fca9549ff81d630b6a91ca0.de-de.xlf"
$dede.Range("J3").Value = "291c5604-b16a-4a9d-9f37-928989d8f104.47401c9d65c9b1103fca9549ff81d630b6a91ca0.de-de.xlf"
$dede.Range("K2").Value = "2016-08-24 15:14:18"
$dede.Range("K3").Value = "2016-08-24 15:14:18"

# ---------------------------------------------------------------
# 4. Column I ("Latest Target File") becomes a hyperlink to the
#    source doc for both language sheets, rows 2 and 3.
# ---------------------------------------------------------------
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $srcDocUrl, "", "", $srcDocName)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $srcDocUrl, "", "", $srcDocName)

$dede.Hyperlinks.Add($dede.Range("I2"), $srcDocUrl, "", "", $srcDocName)
$dede.Hyperlinks.Add($dede.Range("I3"), $srcDocUrl, "", "", $srcDocName)

# ---------------------------------------------------------------
# 5. Widen columns so the longer handback values/links are visible.
# ---------------------------------------------------------------
$ovw.Range("E1").ColumnWidth = 29.9777047293527
$ovw.Range("F1").ColumnWidth = 29.9777047293527

$zhcn.Range("C1").ColumnWidth = 29.9777047293527
$zhcn.Range("I1").ColumnWidth = 40
$zhcn.Range("J1").ColumnWidth = 40

$dede.Range("C1").ColumnWidth = 29.9777047293527
$dede.Range("I1").ColumnWidth = 40
$dede.Range("J1").ColumnWidth = 40

Write-Output "Handback report generated"
